# Add the team's season record (Wins / Losses / Ties) as three new
# columns (AD, AE, AF) to the right of the existing player-stats table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: clone the formatting of an existing header cell (bold,
# centered, bordered) onto the three new header cells, then set their text.
$ws.Range("A1").Copy($ws.Range("AD1"))
$ws.Range("A1").Copy($ws.Range("AE1"))
$ws.Range("A1").Copy($ws.Range("AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-49) gets the same team season record: 73 wins,
# 89 losses, 0 ties.
for ($row = 2; $row -le 49; $row++) {
    $ws.Range("AD" + $row).Value = 73
    $ws.Range("AE" + $row).Value = 89
    $ws.Range("AF" + $row).Value = 0
}

Write-Output "Applied season record columns AD:AF for rows 1-49"
